# feat: add 2022-Q3 data
#
# - Insert a new "2022-Q3" sheet (a copy of "2022-Q2", updated with the new
#   quarter's figures) right after "总计" and before "2022-Q2".
# - Add the corresponding summary row on "总计", re-numbering the existing
#   index column as the old rows shift down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q3" sheet by copying "2022-Q2" (keeps styles,
#    text-vs-number cell types, headers, etc. identical) and placing the
#    copy right after "总计".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$total = $wb.Worksheets.Item("总计")
$q2.Copy($null, $total)

$q3 = $wb.ActiveSheet
$q3.Name = "2022-Q3"

# Update the fund figures on the new "2022-Q3" sheet. The numeric-looking
# columns (D:G) are stored as text in this workbook, same as the sheet they
# were copied from, so force a text format before writing the new values.
$q3.Range("C2").Value = "融通成长30灵活配置混合A/B"

$q3.Range("D2:G3").NumberFormat = "@"

$q3.Range("D2").Value = "1.56"
$q3.Range("E2").Value = "93.13"
$q3.Range("F2").Value = "4.17"
$q3.Range("G2").Value = "0.0651"
$q3.Range("H2").Value = 8

$q3.Range("D3").Value = "1.34"
$q3.Range("E3").Value = "93.13"
$q3.Range("F3").Value = "4.17"
$q3.Range("G3").Value = "0.0559"
$q3.Range("H3").Value = 8

# ---------------------------------------------------------------------
# 2) Add the "2022-Q3" row to "总计", pushing the existing rows down and
#    renumbering the leading index column.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# New row 2 inherits the formatting of what is now row 3 (the old row 2);
# B2:D2 should stay unstyled, like the cells they are modeled after.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.12

# Renumber the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
